# Week 17 data logging for Buccaneers Players Data workbook
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Rushing sheet: update R.Jones / L.Bell / M.Evans rushing stats
# ---------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Row 3 - R.Jones
$rushing.Cells.Item(3, 3).Value = 145
$rushing.Cells.Item(3, 4).Value = 44
$rushing.Cells.Item(3, 5).Value = 21
$rushing.Cells.Item(3, 6).Value = 44

# Row 5 - L.Bell
$rushing.Cells.Item(5, 3).Value = 56
$rushing.Cells.Item(5, 4).Value = 22
$rushing.Cells.Item(5, 6).Value = 14

# Row 6 - M.Evans
$rushing.Cells.Item(6, 3).Value = 4
$rushing.Cells.Item(6, 6).Value = 1

# Leave the Rushing sheet selection parked on C27, matching the saved file
$rushing.Range("C27").Select()

# ---------------------------------------------------------------
# Receiving sheet: A.Brown released (row removed), L.Bell & M.Evans
# added as new receivers, and every other player's season totals
# bumped by their Week 17 production. C.Brate becomes a new row 14.
# ---------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Copy the bold/bordered rank-column style down onto the brand new row 14
$receiving.Range("A13").Copy()
$receiving.Range("A14").PasteSpecial(-4122)

$receivingRows = @(
    @(2,  0,  "R.Jones",      34, 31,  2,  1, 11,  9),
    @(3,  1,  "G.Bernard",    26, 21,  2,  1,  7,  4),
    @(4,  2,  "K.Vaughn",      6,  2,  0,  0,  0,  0),
    @(5,  3,  "L.Bell",        3,  3,  0,  0,  0,  0),
    @(6,  4,  "M.Evans",      73, 51, 34, 18, 17, 12),
    @(7,  6,  "S.Miller",     20, 15,  6,  4,  2,  1),
    @(8,  7,  "T.Johnson",    39, 26,  8,  4,  6,  2),
    @(9,  8,  "J.Darden",      7,  5,  3,  1,  0,  0),
    @(10, 9,  "C.Grayson",     7,  6,  5,  4,  1,  1),
    @(11, 10, "B.Perriman",   10,  4,  3,  2,  2,  1),
    @(12, 11, "R.Gronkowski", 59, 40, 27, 17, 12,  7),
    @(13, 12, "O.Howard",     18, 13,  3,  1,  2,  1),
    @(14, 13, "C.Brate",      47, 27,  6,  0, 19, 10)
)

foreach ($row in $receivingRows) {
    $r = $row[0]
    $receiving.Cells.Item($r, 1).Value = $row[1]
    $receiving.Cells.Item($r, 2).Value = $row[2]
    $receiving.Cells.Item($r, 3).Value = $row[3]
    $receiving.Cells.Item($r, 4).Value = $row[4]
    $receiving.Cells.Item($r, 5).Value = $row[5]
    $receiving.Cells.Item($r, 6).Value = $row[6]
    $receiving.Cells.Item($r, 7).Value = $row[7]
    $receiving.Cells.Item($r, 8).Value = $row[8]
}

# Receiving stays the active sheet/tab, with its selection moved to J11
$receiving.Activate()
$receiving.Range("J11").Select()
